# Ministry spreadsheet update: "overseasAndContact" source-of-infection
# category is being split into "investigating" (still under investigation)
# and a new "established" (source established) category, and three more
# days of daily case data are appended (rows for 2020-03-31 .. 2020-04-02).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates -----------------------------------------------
# Q1 used to be "overseasAndContact"; it is now "investigating".
$ws.Range("Q1").Value = "investigating"
# A brand new "established" column is appended after "community" (R).
$ws.Range("S1").Value = "established"

# --- Row 33 no longer carries the (now retired) O:R breakdown ---------
# Those values move to the new, wider breakdown that starts at row 34.
$ws.Range("O33:R33").ClearContents()

# --- New daily rows 34-36 ----------------------------------------------
# Columns: A date, B confirmed, C totalConfirmed, D probable,
# E totalProbable, F total, G cumulative, H recovered, I totalRecovered,
# J inHospitalNow, K totalBeenInHospital, L inIcu, M deaths, N totalDeaths,
# O overseas, P contact, Q investigating, R community, S established
$newRows = @(
    @{ Row = 34; A = 43921; B = 48; C = 600; D = 10; E = 47; F = 58; G = 647; H = 11; I = 74;  J = 14; L = 2; M = 0; N = 1; O = 343; P = 188; Q = 110; R = 6; S = 647 },
    @{ Row = 35; A = 43922; B = 47; C = 647; D = 14; E = 61; F = 61; G = 708; H = 9;  I = 83;  J = 14; L = 2; M = 0; N = 1; O = 361; P = 212; Q = 127; R = 7; S = 708 },
    @{ Row = 36; A = 43923; B = 76; C = 723; D = 13; E = 74; F = 89; G = 797; H = 9;  I = 92;  J = 13; L = 2; M = 0; N = 1; O = 406; P = 247; Q = 135; R = 8; S = 797 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    foreach ($col in @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "L", "M", "N", "O", "P", "Q", "R", "S")) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
}

# Date column keeps the same custom date display format as the rows above.
$ws.Range("A34:A36").NumberFormat = $ws.Range("A33").NumberFormat

# --- Dimension -----------------------------------------------------------
# Excel keeps the worksheet dimension in sync automatically as cells are
# written (A1:R33 -> A1:S36), so no explicit action is required here.
